$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 35721230
$ws.Range("I70").Value = 3406.75
$ws.Range("J70").Value = 83345000
$ws.Range("K70").Value = 10220.25
$ws.Range("L70").Value = 250035000
$ws.Range("M70").Value = -9950.25
$ws.Range("N70").Value = -250035540
$ws.Range("H73").Value = 35721230
$ws.Range("I73").Value = 3406.75
$ws.Range("J73").Value = 83345000
$ws.Range("K73").Value = 10220.25
$ws.Range("L73").Value = 250035000
$ws.Range("M73").Value = -9284.25
$ws.Range("N73").Value = -250036872
$ws.Range("H76").Value = 7048.9
$ws.Range("I76").Value = 6757.6
$ws.Range("K76").Value = 6757.6
$ws.Range("M76").Value = -6442.6
$ws.Range("H79").Value = 7048.9
$ws.Range("I79").Value = 6757.6
$ws.Range("K79").Value = 6757.6
$ws.Range("M79").Value = -5665.6
$ws.Range("H80").Value = 47628420
$ws.Range("I80").Value = 111112140
$ws.Range("J80").Value = 15622.583
$ws.Range("K80").Value = 333336420
$ws.Range("L80").Value = 46867.749
$ws.Range("M80").Value = -333335422
$ws.Range("N80").Value = -48863.749
$ws.Range("H83").Value = 47628420
$ws.Range("I83").Value = 111112140
$ws.Range("J83").Value = 15622.583
$ws.Range("K83").Value = 1000009260
$ws.Range("L83").Value = 140603.247
$ws.Range("M83").Value = -1000004268
$ws.Range("N83").Value = -150587.247
$ws.Range("H109").Value = 113000
$ws.Range("J109").Value = 113000
$ws.Range("L109").Value = 113000
$ws.Range("N109").Value = -115774

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 22749.75
$ws.Range("I43").Value = 12999
$ws.Range("J43").Value = 26000
$ws.Range("K43").Value = 12999
$ws.Range("L43").Value = 26000
$ws.Range("M43").Value = -12686
$ws.Range("N43").Value = -26626
$ws.Range("H45").Value = 6454.385
$ws.Range("I45").Value = 6434.1113
$ws.Range("K45").Value = 6434.1113
$ws.Range("M45").Value = -6057.1113
$ws.Range("H110").Value = 1692.3334
$ws.Range("I110").Value = 1640.8
$ws.Range("J110").Value = 1950
$ws.Range("K110").Value = 1640.8
$ws.Range("L110").Value = 1950
$ws.Range("M110").Value = 404.2
$ws.Range("N110").Value = -6040

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1766.9048
$ws.Range("I86").Value = 1667.6154
$ws.Range("K86").Value = 1667.6154
$ws.Range("M86").Value = -544.6153999999999
$ws.Range("H89").Value = 1766.9048
$ws.Range("I89").Value = 1667.6154
$ws.Range("K89").Value = 8338.076999999999
$ws.Range("M89").Value = -2722.076999999999
$ws.Range("H128").Value = 6714.9
$ws.Range("I128").Value = 6714.9
$ws.Range("K128").Value = 20144.7
$ws.Range("M128").Value = -17654.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1233.8334
$ws.Range("I16").Value = 1233.8334
$ws.Range("K16").Value = 1233.8334
$ws.Range("M16").Value = -946.8334
$ws.Range("H31").Value = 5121.921
$ws.Range("I31").Value = 2265.111
$ws.Range("J31").Value = 6008.517
$ws.Range("K31").Value = 2265.111
$ws.Range("L31").Value = 6008.517
$ws.Range("M31").Value = -1970.111
$ws.Range("N31").Value = -6598.517
$ws.Range("H34").Value = 5121.921
$ws.Range("I34").Value = 2265.111
$ws.Range("J34").Value = 6008.517
$ws.Range("K34").Value = 2265.111
$ws.Range("L34").Value = 6008.517
$ws.Range("M34").Value = -2063.111
$ws.Range("N34").Value = -6412.517
$ws.Range("H113").Value = 1233.8334
$ws.Range("I113").Value = 1233.8334
$ws.Range("K113").Value = 1233.8334
$ws.Range("M113").Value = 936.1666
$ws.Range("H122").Value = 4426.3784
$ws.Range("I122").Value = 3334.7917
$ws.Range("J122").Value = 6441.615
$ws.Range("K122").Value = 10004.3751
$ws.Range("L122").Value = 19324.845
$ws.Range("M122").Value = -7554.375100000001
$ws.Range("N122").Value = -24224.845
$ws.Range("H132").Value = 3957.0417
$ws.Range("I132").Value = 3638.889
$ws.Range("J132").Value = 4911.5
$ws.Range("K132").Value = 10916.667
$ws.Range("L132").Value = 14734.5
$ws.Range("M132").Value = -8386.667000000001
$ws.Range("N132").Value = -19794.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 418911.22
$ws.Range("I132").Value = 1657.5834
$ws.Range("J132").Value = 836164.8
$ws.Range("K132").Value = 14918.2506
$ws.Range("L132").Value = 7525483.2
$ws.Range("M132").Value = -12388.2506
$ws.Range("N132").Value = -7530543.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3765.25
$ws.Range("I80").Value = 3676
$ws.Range("K80").Value = 3676
$ws.Range("M80").Value = -2678
$ws.Range("H83").Value = 3765.25
$ws.Range("I83").Value = 3676
$ws.Range("K83").Value = 18380
$ws.Range("M83").Value = -13388
$ws.Range("H97").Value = 1786.2727
$ws.Range("I97").Value = 842.2857
$ws.Range("K97").Value = 842.2857
$ws.Range("M97").Value = -346.2857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 355.36365
$ws.Range("I55").Value = 305
$ws.Range("K55").Value = 305
$ws.Range("M55").Value = -132
$ws.Range("H68").Value = 9000
$ws.Range("I68").Value = 20000
$ws.Range("J68").Value = 3500
$ws.Range("K68").Value = 20000
$ws.Range("L68").Value = 3500
$ws.Range("M68").Value = -19251
$ws.Range("N68").Value = -4998
$ws.Range("H71").Value = 9000
$ws.Range("I71").Value = 20000
$ws.Range("J71").Value = 3500
$ws.Range("K71").Value = 100000
$ws.Range("L71").Value = 17500
$ws.Range("M71").Value = -96256
$ws.Range("N71").Value = -24988
$ws.Range("H107").Value = 7745
$ws.Range("I107").Value = 7745
$ws.Range("K107").Value = 7745
$ws.Range("M107").Value = -5825
$ws.Range("H119").Value = 99994.5
$ws.Range("J119").Value = 99994.5
$ws.Range("L119").Value = 99994.5
$ws.Range("N119").Value = -109670.5
$ws.Range("H136").Value = 10431.228
$ws.Range("I136").Value = 4857.6665
$ws.Range("J136").Value = 17119.5
$ws.Range("K136").Value = 14572.9995
$ws.Range("L136").Value = 51358.5
$ws.Range("M136").Value = -12022.9995
$ws.Range("N136").Value = -56458.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5073.778
$ws.Range("I81").Value = 4133
$ws.Range("K81").Value = 8266
$ws.Range("M81").Value = -7205
$ws.Range("H84").Value = 5073.778
$ws.Range("I84").Value = 4133
$ws.Range("K84").Value = 41330
$ws.Range("M84").Value = -36026
$ws.Range("H122").Value = 5552
$ws.Range("J122").Value = 4177.25
$ws.Range("L122").Value = 12531.75
$ws.Range("N122").Value = -17431.75
$ws.Range("H123").Value = 90429
$ws.Range("J123").Value = 90429
$ws.Range("L123").Value = 90429
$ws.Range("N123").Value = -100229
$ws.Range("H132").Value = 2263.5588
$ws.Range("I132").Value = 2257.0386
$ws.Range("J132").Value = 2284.75
$ws.Range("K132").Value = 6771.1158
$ws.Range("L132").Value = 6854.25
$ws.Range("M132").Value = -4241.1158
$ws.Range("N132").Value = -11914.25
$ws.Range("H136").Value = 11826218
$ws.Range("I136").Value = 1685.5
$ws.Range("J136").Value = 18275964
$ws.Range("K136").Value = 5056.5
$ws.Range("L136").Value = 54827892
$ws.Range("M136").Value = -2506.5
$ws.Range("N136").Value = -54832992
